$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "50.875.88"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.899.44"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "366.44"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.84"
$ws.Range("E6").Value = "  -5.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("E7").Value = "  -2.61%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  -4.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.18"
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0827"
$ws.Range("E12").Value = "  -2.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.17"
$ws.Range("E13").Value = "  -3.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.350.15"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.31"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.892.19"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.916"
$ws.Range("E17").Value = "  -4.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "50.829.00"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.20"
$ws.Range("E19").Value = "  -6.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.12"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.72"
$ws.Range("E21").Value = "  -5.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0937"
$ws.Range("E22").Value = "  -3.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.75"
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "256.87"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.65"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.33"
$ws.Range("E26").Value = "  +2.61%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.41"
$ws.Range("E29").Value = "  -3.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.93"
$ws.Range("E30").Value = "  -5.95%  "
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.81"
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.79"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "33.81"
$ws.Range("E36").Value = "  -5.27%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0418"
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("E39").Value = "  -5.53%  "
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.86"
$ws.Range("E41").Value = "  -4.36%  "
$ws.Range("E42").Value = "  -5.92%  "
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.61"
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.63"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.007.71"
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("E49").Value = "  -5.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.183.67"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("E51").Value = "  -2.05%  "
